$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their values as text (many look numeric,
# e.g. "3.44" or "67.044.12", and Excel would otherwise silently convert
# them to numbers/dates when assigned through .Value).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.044.12"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "3.826.42"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "447.60"
$ws.Range("E5").Value = "  +6.62%  "
$ws.Range("D6").Value = "147.57"
$ws.Range("E6").Value = "  +14.71%  "
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  +4.40%  "
$ws.Range("D8").Value = "0.998"
$ws.Range("D9").Value = "0.741"
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("D11").Value = "0.0000326"
$ws.Range("E11").Value = "  -5.52%  "
$ws.Range("D12").Value = "43.67"
$ws.Range("E12").Value = "  +9.87%  "
$ws.Range("D13").Value = "10.36"
$ws.Range("E13").Value = "  +3.60%  "
$ws.Range("D14").Value = "4.418.80"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "15.22"
$ws.Range("E15").Value = "  -4.10%  "
$ws.Range("D16").Value = "3.803.24"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "19.96"
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("D19").Value = "1.15"
$ws.Range("E19").Value = "  +7.52%  "
$ws.Range("D20").Value = "67.009.08"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "424.42"
$ws.Range("E21").Value = "  +4.91%  "
$ws.Range("D22").Value = "14.68"
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("E23").Value = "  +8.67%  "
$ws.Range("D24").Value = "86.60"
$ws.Range("E24").Value = "  +3.93%  "
$ws.Range("D25").Value = "37.45"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("E26").Value = "  +8.07%  "
$ws.Range("D27").Value = "5.49"
$ws.Range("E27").Value = "  -3.82%  "
$ws.Range("D28").Value = "9.50"
$ws.Range("E28").Value = "  +16.65%  "
$ws.Range("D29").Value = "9.72"
$ws.Range("E29").Value = "  +4.73%  "
$ws.Range("D30").Value = "747.23"
$ws.Range("E30").Value = "  +6.82%  "
$ws.Range("D31").Value = "13.73"
$ws.Range("E31").Value = "  +12.45%  "
$ws.Range("D32").Value = "0.134"
$ws.Range("E32").Value = "  +12.22%  "
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").Value = "42.94"
$ws.Range("E34").Value = "  +13.58%  "
$ws.Range("D35").Value = "0.156"
$ws.Range("E35").Value = "  +4.21%  "
$ws.Range("D36").Value = "57.81"
$ws.Range("E36").Value = "  +6.11%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "5.52"
$ws.Range("E38").Value = "  +18.83%  "
$ws.Range("D39").Value = "0.0475"
$ws.Range("E39").Value = "  +5.94%  "
$ws.Range("D40").Value = "0.345"
$ws.Range("E40").Value = "  +16.95%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0683"
$ws.Range("E41").Value = "  -10.96%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "2.87"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +5.04%  "
$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D45").Value = "3.44"
$ws.Range("E45").Value = "  +3.90%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.23"
$ws.Range("E46").Value = "  +5.25%  "
$ws.Range("E47").Value = "  +13.26%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "2.12"
$ws.Range("E48").Value = "  +4.92%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "146.64"
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").Value = "2.64"
$ws.Range("E50").Value = "  +5.67%  "
$ws.Range("D51").Value = "2.87"
$ws.Range("E51").Value = "  +4.84%  "
Write-Host "Applied crypto price/volume update."
